# Updated testing for type 5
# Rebuilds the "Test results" sheet with the full Type 1-4 test matrix and
# adds the per-payment "paid" flag (column E) plus a stray date value on
# the "Payments" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "Test results" sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Test results")

# Header row (Accuracy column moved from old junk index -> just re-assert
# the header text so the shared-string table matches).
$ws.Cells.Item(1,1).Value = "Workflow type"
$ws.Cells.Item(1,2).Value = "Instance"
$ws.Cells.Item(1,3).Value = "Model"
$ws.Cells.Item(1,4).Value = "Time"
$ws.Cells.Item(1,5).Value = "Accuracy"

# Clear out the old partial data below the header before rewriting it.
$ws.Range("A2:E7").Clear() | Out-Null

function Set-TestRow {
    param($Row, $Workflow, $Instance, $Model, $Time, $Accuracy)

    $ws.Cells.Item($Row,1).Value = $Workflow
    if ($null -ne $Instance) {
        $ws.Cells.Item($Row,2).Value = $Instance
    }
    $ws.Cells.Item($Row,2).NumberFormat = "mm-dd-yy"
    $ws.Cells.Item($Row,3).Value = $Model
    $ws.Cells.Item($Row,4).Value = $Time
    $ws.Cells.Item($Row,5).Value = $Accuracy
    $ws.Cells.Item($Row,5).NumberFormat = "0%"
}

Set-TestRow 2  "Type 4" 46327 "gpt-oss:20b"      "525m" 0.89
Set-TestRow 3  "Type 4" 46034 "deepseek-r1:14b"  "192m" 0.78
Set-TestRow 4  "Type 4" 46033 "Qwen3:8b"         "220m" 0.67
Set-TestRow 5  "Type 4" 46033 "llama3.1:8bn"     "183m" 0

Set-TestRow 6  "Type 3" $null "llama3.1:8bn"     "2m"   0
Set-TestRow 7  "Type 3" $null "qwen3:8bn"        "11m"  0.11
Set-TestRow 8  "Type 3" $null "deepseek-r1:14b"  "11m"  0.28000000000000003
Set-TestRow 9  "Type 3" $null "gpt-oss:20b"      "50m"  0.78

Set-TestRow 10 "Type 2" $null "llama3.1:8bn"     "3m"   0
Set-TestRow 11 "Type 2" $null "qwen3:8bn"        "9m"   0
Set-TestRow 12 "Type 2" $null "deepseek-r1:14b"  "NA"   0
Set-TestRow 13 "Type 2" $null "gpt-oss:20b"      "50m"  0

Set-TestRow 14 "Type 1" 46027 "llama3.1:8bn"     "5m"   0
Set-TestRow 15 "Type 1" 46027 "qwen3:8bn"        "83m"  0
Set-TestRow 16 "Type 1" 46027 "deepseek-r1:14b"  "32m"  0
Set-TestRow 17 "Type 1" 46027 "gpt-oss:20b"      "45m"  0.56000000000000005

# Column B ("Instance") holds the raw test-run date but is no longer shown.
$ws.Columns.Item(2).Hidden = $true

# Turn on the header filter over the full table and register the hidden
# _FilterDatabase name Excel normally creates alongside it.
$ws.Range("A1:E17").AutoFilter() | Out-Null
$fdName = $ws.Names.Add("_xlnm._FilterDatabase", "='Test results'!`$A`$1:`$E`$17")
$fdName.Visible = $false

$ws.Range("E27").Select() | Out-Null

# ---------------------------------------------------------------------
# "Payments" sheet
# ---------------------------------------------------------------------
$wsPay = $wb.Worksheets.Item("Payments")

foreach ($r in 2..8) {
    $wsPay.Cells.Item($r,5).Value = 1
}
$wsPay.Cells.Item(10,5).Value = 1

$wsPay.Cells.Item(12,5).Value = 46243
$wsPay.Cells.Item(12,5).NumberFormat = "d-mmm"

$wsPay.Range("E12").Select() | Out-Null
